$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'293.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-5.06%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.12%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.026"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.73%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07381"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.05%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.306"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.29%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.546"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-5.96%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9238"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.86%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1188"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-4.70%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1755"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-3.81%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.20%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-1.45%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.33%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001272"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.88%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.93%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.377"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.84%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'0.3296"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.19%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.587"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.07%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-3.37%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2806"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.47%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03816"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-5.18%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001283"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.32%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003904"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-4.54%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001293"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.62%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003730"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-95.03%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02311"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-9.39%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.55%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007728"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.60%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'CEJI"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.004309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'131.53%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'BKEXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.1275"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.89%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007392"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'10.57%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.006971"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-13.52%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3189"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.06%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006462"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.19%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-8.21%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'35.63%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.05%"
$ws.Range("E51").Style = "Normal"
